$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data was previously read starting one row too low, leaving an empty
# row 1 above the header and a stray leftover row (with only placeholder
# "x" values in C/E) at the bottom. Completing the "read excel file" fix
# means shifting the whole table up by one row: delete the blank row 1
# (rows 2-8 shift up to 1-7), then delete the now-empty leftover row 7
# (the old row 8) that only held the "x" placeholders.
$ws.Rows("1").Delete()
$ws.Rows("7").Delete()

# Restore the cursor position left behind after the cleanup.
$ws.Range("C9").Select()
